# DMS: Translate AppUserStoreMapping Export and ExportTemplate
# Rename the worksheet "AppUserStore" to its Vietnamese translation "Phạm vi đi tuyến"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Name = "Phạm vi đi tuyến"
